$d = $word.ActiveDocument

# --- Title ---------------------------------------------------------------
$d.Content.Find.Execute(
    "Unraveling the Mysteries of Nanotechnology", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "Unveiling the Enigmatic Symphony of Politics in High School", 2) | Out-Null

# --- Byline (merges "Dr" + "." + " Amelia Williams" into one run) --------
$d.Content.Find.Execute(
    "Dr. Amelia Williams", $true, $false, $false, $false, $false, $true, 1,
    $false, "Evelyn Parker", 2) | Out-Null

# --- Email address ---------------------------------------------------------
$d.Content.Find.Execute(
    "awilliams@nanoinstitute", $true, $false, $false, $false, $false, $true,
    1, $false, "evelyn", 2) | Out-Null
$d.Content.Find.Execute(
    "org", $true, $false, $false, $false, $false, $true, 1, $false,
    "parker@schoolmail", 2) | Out-Null

$emailRng = $d.Content
$emailRng.Find.Execute("parker@schoolmail") | Out-Null
$emailEnd = $emailRng.Duplicate
$emailEnd.Collapse(0)
$emailEnd.InsertAfter(".com")

# --- Body paragraph, first block (before the first double line-break) ----
$d.Content.Find.Execute(
    "In the vast realm of scientific exploration, nanotechnology emerges as a transformative force, unveiling a universe of possibilities at the atomic and molecular level",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Politics, the enigmatic symphony of power, influence, and decision-making, plays a pivotal role in shaping our world", 2) | Out-Null

$d.Content.Find.Execute(
    " This interdisciplinary field wields the power to manipulate matter at a scale so minute that it defies comprehension, holding immense promise for breakthroughs across diverse industries",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " From local governance to international relations, politics governs how societies are structured and how individuals interact within them", 2) | Out-Null

$d.Content.Find.Execute(
    " From medicine to manufacturing, and from energy to computing, nanotechnology stands poised to redefine the very fabric of our world. As we delve deeper into its intricate tapestry, we uncover a myriad of applications that hold the potential to revolutionize our technological capabilities and reshape the course of human history",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " As high school students embarking on their journey into adulthood, it is essential to delve into the multifaceted tapestry of politics, unraveling the enigma that surrounds it", 2) | Out-Null

# --- Body paragraph, second block -----------------------------------------
$d.Content.Find.Execute(
    "Nanotechnology's transformative potential lies in its ability to manipulate materials and structures at the nanoscale, where unique properties emerge that are distinct from those observed at larger scales",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "We live in a society governed by laws, policies, and regulations that impact our lives in myriad ways", 2) | Out-Null

$d.Content.Find.Execute(
    " By harnessing these unique properties, scientists have developed groundbreaking technologies that are transforming industries and redefining the limits of what is possible",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " Understanding how these rules are made, who makes them, and why they matter empowers us to engage meaningfully in civic discourse and to hold our elected officials accountable", 2) | Out-Null

$d.Content.Find.Execute(
    " From self-cleaning surfaces and targeted drug delivery systems to ultra-efficient solar cells and lightweight, high-strength materials, the impact of nanotechnology is already being felt across diverse sectors, promising to usher in a new era of innovation and progress",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " Politics is not just about abstract theories and power struggles; it is about real people, their needs, and their aspirations", 2) | Out-Null

# --- Body paragraph, third block ------------------------------------------
$d.Content.Find.Execute(
    "As research continues to unravel the mysteries of nanotechnology, we are witnessing an acceleration in the development of cutting-edge applications that hold immense promise for addressing some of the world's most pressing challenges",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Moreover, politics is deeply intertwined with history, economics, and culture", 2) | Out-Null

$d.Content.Find.Execute(
    " From developing new cancer treatments and targeted drug delivery systems to harnessing nanomaterials for sustainable energy solutions, nanotechnology is poised to make significant contributions to human health and well-being, environmental sustainability, and economic prosperity",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " By examining the interplay of these factors, we can gain valuable insights into the complexities of the political landscape", 2) | Out-Null

$d.Content.Find.Execute(
    " This emerging field is a testament to the boundless potential of human ingenuity and serves as a beacon of hope for a brighter, more technologically advanced future",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " As we navigate the ever-changing political climate, we must cultivate critical thinking skills, learn to analyze information objectively, and develop a deep appreciation for diverse perspectives", 2) | Out-Null

# --- Summary paragraph (also drops the lastRenderedPageBreak hint) -------
$d.Content.Find.Execute(
    "Nanotechnology, a rapidly evolving field that explores the manipulation of matter at the atomic and molecular scale, holds immense promise for transformative applications across diverse industries",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Politics is an intricate and dynamic field that encompasses power, influence, decision-making, and governance", 2) | Out-Null

$d.Content.Find.Execute(
    " With its unique ability to manipulate materials and structures at the nanoscale, nanotechnology is driving the development of groundbreaking technologies that are redefining the limits of what is possible",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " It plays a vital role in shaping societies, affecting the lives of individuals in numerous ways", 2) | Out-Null

$d.Content.Find.Execute(
    " From self-cleaning surfaces and targeted drug delivery systems to ultra-efficient solar cells and lightweight, high-strength materials, the impact of nanotechnology is already being felt across various sectors",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " By studying politics, high school students can cultivate critical thinking skills, appreciate diverse perspectives, and gain valuable insights into the world around them", 2) | Out-Null

$d.Content.Find.Execute(
    " As research continues to unveil the mysteries of nanotechnology, we can anticipate further breakthroughs that will address global challenges, improve human health and well-being, and usher in a new era of scientific advancement",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " As active citizens, they can participate meaningfully in civic discourse, hold leaders accountable, and contribute to a more just and equitable society", 2) | Out-Null

# New trailing sentence appended after the summary's final sentence
$sumRng = $d.Content
$sumRng.Find.Execute("a more just and equitable society") | Out-Null
$sumEnd = $sumRng.Duplicate
$sumEnd.Collapse(0)
$sumEnd.InsertAfter(". Understanding politics empowers us to navigate the complexities of the modern world and to make informed decisions that shape our future")

# --- New trailing empty paragraph at the very end of the document --------
$d.Paragraphs.Last.Range.InsertParagraphAfter()

Write-Host "edit applied"
